$d = $word.ActiveDocument

# 1. Remove the stray <w:lastRenderedPageBreak/> marker that precedes the
#    "Alternative flow of events" heading in the first use case. Scope the
#    Find/Replace to that single paragraph's range so only that run is
#    regenerated (and the other two "Alternative flow of events" headings
#    elsewhere in the document are left untouched).
$targetText = "Alternative flow of events"
$heading2Count = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text -eq $targetText) {
        $heading2Count = $heading2Count + 1
        if ($heading2Count -eq 1) {
            $pr = $p.Range
            $pr.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, $targetText, 2)
        }
    }
}

# 2. Every body paragraph that is NOT one of the heading/title styles gets
#    "space after paragraph" collapsed to 0 (i.e. <w:spacing w:after="0"/>
#    inserted into its <w:pPr>).
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -ne "Title" -and $styleName -ne "Heading 1" -and $styleName -ne "Heading 2" -and $styleName -ne "Heading 3") {
        $p.SpaceAfter = 0
    }
}
